$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

$ws.Range("A27").Value = "vaishali.kh2310@gmail.com"
$ws.Range("B27").Value = "Login"
$ws.Range("C27").Value = "2025-06-17 21:53:31"

$ws.Range("A28").Value = "vaishali.kh2310@gmail.com"
$ws.Range("B28").Value = "Login"
$ws.Range("C28").Value = "2025-06-18 11:21:42"
